# Updates the answer table on the worksheet with a newly generated set of
# two-digit-number ÷ one-digit-number problems/answers.
#
# Each populated row of the table (rows 1, 5, 9, 13, 17) holds five
# "dividend÷divisor=quotient, remainder" strings, one per column. We replace
# each cell's text directly by table coordinates, which avoids any ambiguity
# from duplicate/overlapping old & new values elsewhere in the document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellUpdates = @(
    @{ Row = 1;  Col = 1; New = "35÷5=7, 0" },
    @{ Row = 1;  Col = 2; New = "11÷5=2, 1" },
    @{ Row = 1;  Col = 3; New = "22÷6=3, 4" },
    @{ Row = 1;  Col = 4; New = "53÷6=8, 5" },
    @{ Row = 1;  Col = 5; New = "69÷8=8, 5" },

    @{ Row = 5;  Col = 1; New = "87÷2=43, 1" },
    @{ Row = 5;  Col = 2; New = "43÷4=10, 3" },
    @{ Row = 5;  Col = 3; New = "89÷6=14, 5" },
    @{ Row = 5;  Col = 4; New = "14÷6=2, 2" },
    @{ Row = 5;  Col = 5; New = "67÷2=33, 1" },

    @{ Row = 9;  Col = 1; New = "23÷6=3, 5" },
    @{ Row = 9;  Col = 2; New = "79÷7=11, 2" },
    @{ Row = 9;  Col = 3; New = "44÷8=5, 4" },
    @{ Row = 9;  Col = 4; New = "98÷2=49, 0" },
    @{ Row = 9;  Col = 5; New = "74÷8=9, 2" },

    @{ Row = 13; Col = 1; New = "45÷2=22, 1" },
    @{ Row = 13; Col = 2; New = "58÷6=9, 4" },
    @{ Row = 13; Col = 3; New = "43÷2=21, 1" },
    @{ Row = 13; Col = 4; New = "75÷9=8, 3" },
    @{ Row = 13; Col = 5; New = "68÷5=13, 3" },

    @{ Row = 17; Col = 1; New = "10÷2=5, 0" },
    @{ Row = 17; Col = 2; New = "82÷5=16, 2" },
    @{ Row = 17; Col = 3; New = "40÷7=5, 5" },
    @{ Row = 17; Col = 4; New = "97÷7=13, 6" },
    @{ Row = 17; Col = 5; New = "51÷5=10, 1" }
)

foreach ($update in $cellUpdates) {
    $t.Cell($update.Row, $update.Col).Range.Text = $update.New
}
